$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "75÷8="
$tbl.Cell(1, 2).Range.Text = "13÷6="
$tbl.Cell(1, 3).Range.Text = "75÷7="
$tbl.Cell(1, 4).Range.Text = "27÷4="
$tbl.Cell(1, 5).Range.Text = "85÷5="
$tbl.Cell(5, 1).Range.Text = "90÷7="
$tbl.Cell(5, 2).Range.Text = "43÷2="
$tbl.Cell(5, 3).Range.Text = "72÷9="
$tbl.Cell(5, 5).Range.Text = "66÷3="
$tbl.Cell(9, 1).Range.Text = "82÷6="
$tbl.Cell(9, 2).Range.Text = "23÷3="
$tbl.Cell(9, 3).Range.Text = "54÷8="
$tbl.Cell(9, 4).Range.Text = "33÷9="
$tbl.Cell(9, 5).Range.Text = "87÷3="
$tbl.Cell(13, 1).Range.Text = "30÷7="
$tbl.Cell(13, 2).Range.Text = "83÷5="
$tbl.Cell(13, 3).Range.Text = "76÷9="
$tbl.Cell(13, 4).Range.Text = "91÷4="
$tbl.Cell(13, 5).Range.Text = "85÷8="
$tbl.Cell(17, 1).Range.Text = "77÷2="
$tbl.Cell(17, 2).Range.Text = "68÷8="
$tbl.Cell(17, 3).Range.Text = "61÷7="
$tbl.Cell(17, 4).Range.Text = "18÷6="
$tbl.Cell(17, 5).Range.Text = "55÷3="
